$d = $word.ActiveDocument
$W = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Change 1: "Use Case:" paragraph (row 1 of the table)
#   - drop the spellcheck proofErr markers around "Use"
#   - merge "Use" + " Case:" into a single bold run "Use Case:"
#   - append a new run "en einsehen" after " Statistik"
#   - add a _GoBack bookmark at the end of the paragraph
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/></w:rPr><w:t>Use Case:</w:t></w:r><w:r><w:t xml:space="preserve"> Statistik</w:t></w:r><w:r><w:t>en einsehen</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$p1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: "Vorbedingungen:" paragraph - merge the split "Der Spi" / "eler..."
#   runs back into a single run.
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs(8).Range
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Vorbedingungen: </w:t></w:r><w:r><w:t>Der Spieler hat die Anwendung geöffnet und bereits mindestens ein Spiel gespielt (ansonsten macht die Statistik wenig Sinn).</w:t></w:r></w:p>
'@
$p8.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Change 3: "Beschreibung" list-intro paragraph ending in ":" - remove the
#   _GoBack bookmark (it was relocated to the end of the title paragraph).
# ---------------------------------------------------------------------------
$p10 = $d.Paragraphs(10).Range
$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Der Spieler</w:t></w:r><w:r><w:t xml:space="preserve"> kann</w:t></w:r><w:r><w:t xml:space="preserve"> im Hauptmenü den Punkt „Statistik“</w:t></w:r><w:r><w:t xml:space="preserve"> auswählen. I</w:t></w:r><w:r><w:t>m Statistikbildschirm werden alle statistischen Werte angezeigt, die in allen bereits gespielten Spielen erfasst wurden. Darunter fallen</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>
'@
$p10.InsertXML($xml3)

# ---------------------------------------------------------------------------
# Change 4: append a page break and a large new "Beschreibung" section after
#   the table, right before the trailing empty body paragraph.
# ---------------------------------------------------------------------------
$t1 = $d.Tables(1)
$tailRange = $d.Range($t1.Range.End, $t1.Range.End)

$xml4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p/>
<w:p><w:r><w:br w:type="page"/></w:r></w:p>
<w:p><w:r><w:lastRenderedPageBreak/><w:t>Der Spieler w</w:t></w:r><w:r><w:t>ill das bisher von ihm im Spiel Erreichte in Zahlen sehen und mit anderen Spielern vergleichen können. Insbesondere der Highscore ist für ihn als Massstab seiner Leistung interessant, der Rest dient eher als unterhaltsame Trivia.</w:t></w:r></w:p>
<w:p><w:r><w:t>Der Spieler kann im Hauptmenü den Punkt „Statistik“ auswählen. Im Statistikbildschirm werden alle statistischen Werte angezeigt, die in allen bereits gespielten Spielen erfasst wurden. Darunter fallen:</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Highscore für jeden Spielmodus</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Längste Serie im Endlosmodus (Zeitspanne)</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Anzahl gespielter Spiele für jeden Spielmodus</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Anzahl gewonnener Spiele (d.h. das Spiel wurde erfolgreich abgeschlossen)</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Bisherige Spielzeit</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Anzahl verladener Container</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Anzahl beladener Schiffe</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Anzahl zerstörter Schiffe</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Verhältnis von beladenen zu zerstörten Schiffen</w:t></w:r></w:p>
<w:p><w:r><w:t>Die Liste der zu erfassenden Daten ist weder final noch abschliessend, sie wird natürlich den anderen Anforderungen und dem Spielablauf angepasst.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Die Statistik-Funktionen sind niedrig priorisiert. Sollten Werte schwierig zu erfassen sein, dann sollen sie einfach weggelassen werden. Es soll wenig Zeit in die Statistik investiert werden.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$tailRange.InsertXML($xml4)
